# Apply the weekly update described in the commit: two new price records
# (Murcott mandarina, $/bandeja 18 kilos, Región de O'Higgins, 2023-08-28)
# are inserted at the top of the Agro Chillán block, pushing the existing
# rows 378-404 down to 380-406.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows above the current row 378; this also shifts the
# number-format/style of the Fecha column down into the new rows.
$ws.Rows("378:379").Insert()

# --- New row 378 ---
$ws.Cells.Item(378, 1).Value = 7
$ws.Cells.Item(378, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(378, 3).Value = "Ñuble"
$ws.Cells.Item(378, 4).Value = 45166
$ws.Cells.Item(378, 5).Value = 16
$ws.Cells.Item(378, 6).Value = "Fruta"
$ws.Cells.Item(378, 7).Value = 100102
$ws.Cells.Item(378, 8).Value = "Cítricos"
$ws.Cells.Item(378, 9).Value = 100102004
$ws.Cells.Item(378, 10).Value = "Mandarina"
$ws.Cells.Item(378, 11).Value = "Murcott"
$ws.Cells.Item(378, 12).Value = "Primera"
$ws.Cells.Item(378, 13).Value = 100
$ws.Cells.Item(378, 14).Value = 10000
$ws.Cells.Item(378, 15).Value = 10000
$ws.Cells.Item(378, 16).Value = 10000
$ws.Cells.Item(378, 17).Value = "`$/bandeja 18 kilos"
$ws.Cells.Item(378, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(378, 19).Value = 556
$ws.Cells.Item(378, 20).Value = 18

# --- New row 379 ---
$ws.Cells.Item(379, 1).Value = 7
$ws.Cells.Item(379, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(379, 3).Value = "Ñuble"
$ws.Cells.Item(379, 4).Value = 45166
$ws.Cells.Item(379, 5).Value = 16
$ws.Cells.Item(379, 6).Value = "Fruta"
$ws.Cells.Item(379, 7).Value = 100102
$ws.Cells.Item(379, 8).Value = "Cítricos"
$ws.Cells.Item(379, 9).Value = 100102004
$ws.Cells.Item(379, 10).Value = "Mandarina"
$ws.Cells.Item(379, 11).Value = "Murcott"
$ws.Cells.Item(379, 12).Value = "Segunda"
$ws.Cells.Item(379, 13).Value = 120
$ws.Cells.Item(379, 14).Value = 8000
$ws.Cells.Item(379, 15).Value = 8000
$ws.Cells.Item(379, 16).Value = 8000
$ws.Cells.Item(379, 17).Value = "`$/bandeja 18 kilos"
$ws.Cells.Item(379, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(379, 19).Value = 444
$ws.Cells.Item(379, 20).Value = 18
